$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.926.60"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.94%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.263.10"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.38%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "112.54"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.91%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "299.81"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.41%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.635"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.16%  "
$ws.Range("E8").Value = "  -0.33%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.614"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.60%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "44.30"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.20%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0928"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.08%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.62"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.89%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "8.91"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.19%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.07"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +21.02%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.105"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.29%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.39"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.85%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.602.18"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.245.13"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.12%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "42.831.15"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.29%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0000107"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.07%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.23"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.73%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "75.41"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.52%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.54"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +14.68%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "258.35"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +11.54%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.43"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.89%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.96"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.62%  "
$ws.Range("E27").Value = "  -0.20%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.64"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.26%  "
$ws.Range("E29").Value = "  -0.52%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "38.15"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.47%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "22.37"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.98%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "175.13"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.89%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.16"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.41%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0893"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.39%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.71"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.43%  "
$ws.Range("E36").Value = "  +9.24%  "
$ws.Range("E37").Value = "  +0.54%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.24"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -6.08%  "
$ws.Range("E39").Value = "  +1.46%  "
$ws.Range("E40").Value = "  -1.76%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.45"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.63%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "71.84"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.80%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.232"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.86%  "
$ws.Range("E44").Value = "  -0.27%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.55"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -6.71%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.34"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.70%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "5.54"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.06%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "107.33"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +6.87%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.30"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.70%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.71"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.79%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "73.19"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.04%  "
